# Weekly Fruta/Hortaliza update: insert 2 new daily-price rows for "Choclo"
# (Vega Modelo de Temuco) ahead of the existing data block, pushing the
# rest of the table down by two rows (old row 235 -> new row 237, etc.),
# extending the sheet from A1:R317 to A1:R319.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the current row 235 position; everything that
# was at row 235 onward (through the old last row 317) shifts down to
# rows 237 through 319.
$ws.Range("A235:A236").EntireRow.Insert()

# --- New row 235 ---
$ws.Range("A235").Value = 10
$ws.Range("B235").Value = "Vega Modelo de Temuco"
$ws.Range("C235").Value = "La Araucanía"
$ws.Range("D235").Value = 44559
$ws.Range("E235").Value = 9
$ws.Range("F235").Value = 100112024
$ws.Range("G235").Value = "Choclo"
$ws.Range("H235").Value = "Choclero"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 9000
$ws.Range("K235").Value = 380
$ws.Range("L235").Value = 400
$ws.Range("M235").Value = 391
$ws.Range("N235").Value = "`$/unidad"
$ws.Range("O235").Value = "Región del Maule"
$ws.Range("P235").Value = 391
$ws.Range("Q235").Value = 1
$ws.Range("R235").Value = "Hortaliza"

# --- New row 236 ---
$ws.Range("A236").Value = 10
$ws.Range("B236").Value = "Vega Modelo de Temuco"
$ws.Range("C236").Value = "La Araucanía"
$ws.Range("D236").Value = 44559
$ws.Range("E236").Value = 9
$ws.Range("F236").Value = 100112024
$ws.Range("G236").Value = "Choclo"
$ws.Range("H236").Value = "Dulce o Americano"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 12000
$ws.Range("K236").Value = 280
$ws.Range("L236").Value = 300
$ws.Range("M236").Value = 290
$ws.Range("N236").Value = "`$/unidad"
$ws.Range("O236").Value = "Región del Maule"
$ws.Range("P236").Value = 290
$ws.Range("Q236").Value = 1
$ws.Range("R236").Value = "Hortaliza"

# Match the date cell format used by the rest of column D.
$ws.Range("D235:D236").NumberFormat = "YYYY-MM-DD HH:MM:SS"
